$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# YDS sheet: append per-play yardage logs (Rush/Pass, OFF/DEF) with
# the new divisional-round game's numbers.
# ------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 1 2 5 -1 2 2 7 4 2 6 4 5 4 5 1 5 9 5 8 13 5 -2 4 3 0 9 1 1"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 15 16 12 18 3 12 -1 24 6 12 14"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 4 -1 1 5 6 14 3 4 4 3 2 -1 2 3 2 5 4 3 2 2"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 14 19 10 11 7 9 1 3 8 9 5 75 6 6 25 1 5 6 2 4"

# ------------------------------------------------------------------
# OFF sheet: season totals for Home (row 2) and Road (row 3).
# ------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 277
$wsOFF.Range("E2").Value = 14
$wsOFF.Range("F2").Value = 83
$wsOFF.Range("G2").Value = 74
$wsOFF.Range("I2").Value = 12
$wsOFF.Range("J2").Value = 41
$wsOFF.Range("N2").Value = 20
$wsOFF.Range("O2").Value = 21
$wsOFF.Range("B3").Value = 14
$wsOFF.Range("C3").Value = 181
$wsOFF.Range("E3").Value = 41
$wsOFF.Range("F3").Value = 107
$wsOFF.Range("G3").Value = 30
$wsOFF.Range("H3").Value = 36
$wsOFF.Range("I3").Value = 65
$wsOFF.Range("J3").Value = 48
$wsOFF.Range("L3").Value = 313
$wsOFF.Range("M3").Value = 210
$wsOFF.Range("Q3").Value = 681

# ------------------------------------------------------------------
# DEF sheet: season totals for Home (row 2) and Road (row 3).
# ------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 221
$wsDEF.Range("F2").Value = 73
$wsDEF.Range("G2").Value = 69
$wsDEF.Range("J2").Value = 34
$wsDEF.Range("N2").Value = 42
$wsDEF.Range("B3").Value = 15
$wsDEF.Range("C3").Value = 214
$wsDEF.Range("E3").Value = 49
$wsDEF.Range("F3").Value = 116
$wsDEF.Range("H3").Value = 39
$wsDEF.Range("I3").Value = 76
$wsDEF.Range("J3").Value = 71
$wsDEF.Range("L3").Value = 369
$wsDEF.Range("M3").Value = 249
$wsDEF.Range("Q3").Value = 700

# ------------------------------------------------------------------
# ST sheet: season kicking/return totals (row 2), plus per-kick logs
# (RA/RM/D columns, rows 4-6 and row 3) for the new game.
# ------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 96
$wsST.Range("D2").Value = 64
$wsST.Range("F2").Value = 590
$wsST.Range("G2").Value = 575
$wsST.Range("J2").Value = 291
$wsST.Range("K2").Value = 276
$wsST.Range("L2").Value = 164
$wsST.Range("M2").Value = 129
$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " 49 59 61"
$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 46 42 52 49 45"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " 15 26 25"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 0 0 6 5 0"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 32 45 14"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 0 6 0 7 9"

# ------------------------------------------------------------------
# TURNS sheet: Road turnovers count.
# ------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 9

# ------------------------------------------------------------------
# PEN sheet: penalty counts.
# ------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 24
$wsPEN.Range("B3").Value = 23
$wsPEN.Range("D3").Value = 3
